# Rename the existing sheet to "shifts" and add a new sheet "adir" after it,
# then populate both sheets with the shift-scheduling data.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.ActiveSheet
$ws1.Name = "shifts"

$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "adir"

# --- Sheet "shifts" ---
$ws1.Range("B1").Value = "Sunday"
$ws1.Range("C1").Value = "Monday"
$ws1.Range("D1").Value = "Tuesday"
$ws1.Range("E1").Value = "Wednesday"
$ws1.Range("F1").Value = "Thursday"
$ws1.Range("G1").Value = "Friday"
$ws1.Range("H1").Value = "Saturday"

$ws1.Range("A2").Value = "Morning"
$ws1.Range("B2").Value = "adir"
$ws1.Range("C2").Value = "yoni"
$ws1.Range("D2").Value = "tair"
$ws1.Range("E2").Value = "asaf"
$ws1.Range("F2").Value = "adir"
$ws1.Range("G2").Value = "tair"

$ws1.Range("A3").Value = "Morning"
$ws1.Range("B3").Value = "rotem"
$ws1.Range("C3").Value = "adir"
$ws1.Range("D3").Value = "rotem"
$ws1.Range("E3").Value = "stav"
$ws1.Range("F3").Value = "yoni"
$ws1.Range("G3").Value = "stav"

$ws1.Range("A4").Value = "shift manager"
$ws1.Range("B4").Value = "emilia"
$ws1.Range("C4").Value = "michal"
$ws1.Range("D4").Value = "michal"
$ws1.Range("E4").Value = "emilia"
$ws1.Range("F4").Value = "michal"
$ws1.Range("G4").Value = "emilia"

$ws1.Range("A5").Value = "Evening"
$ws1.Range("B5").Value = "asaf"
$ws1.Range("C5").Value = "rotem"
$ws1.Range("D5").Value = "asaf"
$ws1.Range("E5").Value = "yoni"
$ws1.Range("F5").Value = "adir"
$ws1.Range("H5").Value = "rotem"

$ws1.Range("A6").Value = "Evening"
$ws1.Range("B6").Value = "rotem"
$ws1.Range("C6").Value = "tair"
$ws1.Range("D6").Value = "adir"
$ws1.Range("E6").Value = "stav"
$ws1.Range("F6").Value = "yoni"
$ws1.Range("H6").Value = "asaf"

$ws1.Range("A7").Value = "shift manager"
$ws1.Range("B7").Value = "michal"
$ws1.Range("C7").Value = "michal"
$ws1.Range("D7").Value = "emilia"
$ws1.Range("E7").Value = "emilia"
$ws1.Range("F7").Value = "michal"
$ws1.Range("H7").Value = "emilia"

# --- Sheet "adir" ---
$ws2.Range("B1").Value = "Sunday"
$ws2.Range("C1").Value = "Monday"
$ws2.Range("D1").Value = "Tuesday"
$ws2.Range("E1").Value = "Wednesday"
$ws2.Range("F1").Value = "Thursday"
$ws2.Range("G1").Value = "Friday"
$ws2.Range("H1").Value = "Saturday"

$ws2.Range("A2").Value = "Morning"
$ws2.Range("G2").Value = "NO"

$ws2.Range("A3").Value = "Evening"
$ws2.Range("C3").Value = "NO"

$ws1.Activate()
